{"js": "const doc = context.document;\n\n// Create the three new character styles.\ndoc.addStyle(\"GaNStyle\", \"Character\");\ndoc.addStyle(\"GaNParagraph\", \"Character\");\ndoc.addStyle(\"GaNLinks\", \"Character\");\nawait context.sync();\n\n// Re-fetch the styles so the font changes persist correctly.\nconst styles = doc.getStyles();\nconst gaNStyle = styles.getByNameOrNullObject(\"GaNStyle\");\nconst gaNParagraph = styles.getByNameOrNullObject(\"GaNParagraph\");\nconst gaNLinks = styles.getByNameOrNullObject(\"GaNLinks\");\nawait context.sync();\n\ngaNStyle.font.name = \"Calibri\";\ngaNStyle.font.size = 14;\n\ngaNParagraph.font.name = \"Calibri\";\ngaNParagraph.font.size = 10;\n\ngaNLinks.font.name = \"Calibri\";\ngaNLinks.font.bold = true;\ngaNLinks.font.color = \"#000080\";\ngaNLinks.font.size = 9.5;\ngaNLinks.font.underline = \"Single\";\n\nawait context.sync();\n\n// Apply GaNParagraph to every run containing the repeated Swedish paragraph.\nconst body = context.document.body;\nconst paragraphText = \"Du deltar i en v\u00e4rldsomsp\u00e4nnande kampanj f\u00f6r att observera och rapportera de svagaste synliga stj\u00e4rnorna, som ett m\u00e5tt p\u00e5 ljusf\u00f6roreningarna p\u00e5 orten. Genom att hitta och observera Bootes konstellation p\u00e5 natthimlen kan folk i hela v\u00e4rlden l\u00e4ra sig hur belysningen i v\u00e5ra samh\u00e4llen och omgivningar bidrar till ljusf\u00f6roreningar. Era bidrag till online-databasen hj\u00e4lper till att dokumentera den synliga natthimlens \u00f6ver hela v\u00e4rlden.\";\nconst paragraphResults = body.search(paragraphText, { matchCase: true });\nparagraphResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphResults.items.length; i++) {\n  paragraphResults.items[i].style = \"GaNParagraph\";\n}\nawait context.sync();\n\n// Apply GaNLinks to the run with the GaNight map link.\nconst linkText = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\";\nconst linkResults = body.search(linkText, { matchCase: true });\nlinkResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < linkResults.items.length; i++) {\n  linkResults.items[i].style = \"GaNLinks\";\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Add the three new character styles (order matches the target styles.xml)\n$gaNStyle = $d.Styles.Add(\"GaNStyle\", 2)\n$gaNStyle.Font.Name = \"Calibri\"\n$gaNStyle.Font.Size = 14\n\n$gaNParagraph = $d.Styles.Add(\"GaNParagraph\", 2)\n$gaNParagraph.Font.Name = \"Calibri\"\n$gaNParagraph.Font.Size = 10\n\n$gaNLinks = $d.Styles.Add(\"GaNLinks\", 2)\n$gaNLinks.Font.Name = \"Calibri\"\n$gaNLinks.Font.Bold = $true\n$gaNLinks.Font.Color = 8388608\n$gaNLinks.Font.Size = 9.5\n$gaNLinks.Font.Underline = 1\n\n# Apply GaNParagraph to every run containing the repeated Swedish paragraph\n$paragraphText = \"Du deltar i en v\u00e4rldsomsp\u00e4nnande kampanj f\u00f6r att observera och rapportera de svagaste synliga stj\u00e4rnorna, som ett m\u00e5tt p\u00e5 ljusf\u00f6roreningarna p\u00e5 orten. Genom att hitta och observera Bootes konstellation p\u00e5 natthimlen kan folk i hela v\u00e4rlden l\u00e4ra sig hur belysningen i v\u00e5ra samh\u00e4llen och omgivningar bidrar till ljusf\u00f6roreningar. Era bidrag till online-databasen hj\u00e4lper till att dokumentera den synliga natthimlens \u00f6ver hela v\u00e4rlden.\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $paragraphText\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\nwhile ($rng.Find.Execute()) {\n  $rng.Style = \"GaNParagraph\"\n  $rng.Collapse(0)\n}\n\n# Apply GaNLinks to the run with the GaNight map link\n$linkText = \"(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).\"\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = $linkText\n$rng2.Find.Forward = $true\n$rng2.Find.Wrap = 0\nwhile ($rng2.Find.Execute()) {\n  $rng2.Style = \"GaNLinks\"\n  $rng2.Collapse(0)\n}\n"}
